# Update countries & provincias Spain
# Applies the data refresh captured in the commit "Update countries & provincias Spain":
#  - the "Datos actualizados..." timestamp moves from 15:52 to 16:22
#  - Estados Unidos (row 4), Austria (row 19) and Georgia (row 112) get refreshed case counts
#  - Moldavia overtakes Marruecos (rows 61/62), Bulgaria overtakes Uzbekistan (rows 81/82)
#    and Sudan del Sur overtakes Gambia (rows 207/208) in the case-count ranking, so the
#    two countries in each pair swap rows while keeping the table sorted by total cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 10 de Abril de 2020 a las 16:22"

function Set-CountryRow($Row, $Name, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    if ($Name) {
        $ws.Cells.Item($Row, 1).Value = $Name
    }
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Estados Unidos - refreshed counts, no reordering
Set-CountryRow 4 $null 469464 898 25937 426816 10011 20 16711

# Austria - refreshed counts, no reordering
Set-CountryRow 19 $null 13494 250 6064 7111 261 24 319

# Moldavia overtakes Marruecos
Set-CountryRow 61 "Moldavia" 1438 149 56 1353 80 0 29
Set-CountryRow 62 "Marruecos" 1431 57 114 1212 1 8 105

# Bulgaria overtakes Uzbekistan
Set-CountryRow 81 "Bulgaria" 635 17 54 556 33 1 25
Set-CountryRow 82 "Uzbekistan" 624 42 42 579 8 0 3

# Georgia - refreshed counts, no reordering
Set-CountryRow 112 $null 230 12 54 173 6 0 3

# Sudan del Sur overtakes Gambia
Set-CountryRow 207 "Sudan del Sur" 4 1 0 4 0 0 0
Set-CountryRow 208 "Gambia" 4 0 2 1 0 0 1
